$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 25569.00049503472
$ws.Range("C2").Value = "January"
$ws.Range("D2").Value = "Anil"
$ws.Range("E2").Value = "Rajasthan"
$ws.Range("F2").Value = "HWRAJ"
$ws.Range("G2").Value = 2000

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 25569.00049503472
$ws.Range("C3").Value = "January"
$ws.Range("D3").Value = "Anil"
$ws.Range("E3").Value = "Rajasthan"
$ws.Range("F3").Value = "HWRAJ"
$ws.Range("G3").Value = 4000

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 25569.00049503472
$ws.Range("C4").Value = "January"
$ws.Range("D4").Value = "Anil"
$ws.Range("E4").Value = "Rajasthan"
$ws.Range("F4").Value = "HWKTK"
$ws.Range("G4").Value = 4000

# The date column (B) needs a single shared "date" cell style (numFmtId 14,
# i.e. built-in short-date format "mm-dd-yy") applied to B2:B4. Setting
# .NumberFormat on each cell (or on the whole range at once) allocates a
# brand-new style record per cell in this engine, so instead we build the
# style once on a scratch cell, then copy *just the formatting* onto the
# target cells (PasteSpecial xlPasteFormats shares one style index across
# the whole destination), and finally remove the scratch cell again.
$scratch = $ws.Range("Z1")
$scratch.Value = 25569.00049503472
$scratch.NumberFormat = "mm-dd-yy"
$scratch.Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$scratch.Delete(-4159)
